# Emissions ceramics + update notebook residential tertiary
#
# 1. Insert a new worksheet "year_Vecteur" right before "retrofit_Transition"
#    and populate it with Vecteur / year / direct_emissions / indirect_emissions
#    data (5 energy vectors: elec, gaz, fioul, bois, charbon).
# 2. Update the stored cell-selection view-state on a few other sheets.

$wb = $excel.ActiveWorkbook

# --- view-state tweak on "0D" --------------------------------------------
$ws0D = $wb.Worksheets.Item("0D")
$ws0D.Range("A19").Select()

# --- view-state tweak on "Production_system_year" ------------------------
$wsProdYear = $wb.Worksheets.Item("Production_system_year")
$wsProdYear.Range("G1").Select()

# --- new sheet "year_Vecteur", inserted just before "retrofit_Transition" -
$wsRetrofit = $wb.Worksheets.Item("retrofit_Transition")
$wsNew = $wb.Worksheets.Add($wsRetrofit)
$wsNew.Name = "year_Vecteur"

$wsNew.Range("A1").Value = "Vecteur"
$wsNew.Range("B1").Value = "year"
$wsNew.Range("C1").Value = "direct_emissions"
$wsNew.Range("D1").Value = "indirect_emissions"

$wsNew.Range("A2").Value = "elec"
$wsNew.Range("B2").Value = 2020
$wsNew.Range("C2").Value = 0
$wsNew.Range("D2").Value = 79

$wsNew.Range("A3").Value = "gaz"
$wsNew.Range("B3").Value = 2020
$wsNew.Range("C3").Value = 187
$wsNew.Range("D3").Formula = "=227-C3"

$wsNew.Range("A4").Value = "fioul"
$wsNew.Range("B4").Value = 2020
$wsNew.Range("C4").Value = 272
$wsNew.Range("D4").Value = 57

$wsNew.Range("A5").Value = "bois"
$wsNew.Range("B5").Value = 2020
$wsNew.Range("C5").Value = 27
$wsNew.Range("D5").Value = 0

$wsNew.Range("A6").Value = "charbon"
$wsNew.Range("B6").Value = 2020
$wsNew.Range("C6").Value = 346.5
$wsNew.Range("D6").Value = 28.5

# match page setup of the sibling data sheets (A4 portrait)
$wsNew.PageSetup.PaperSize = 9
$wsNew.PageSetup.Orientation = 1

# Leave the new sheet's selection on F13 and make it the active tab (last
# sheet touched/selected ends up being the active one on save).
$wsNew.Range("F13").Select()
